# feat: add 2022-Q1 data
#
# 1. Turn the existing "总计" (summary) sheet into the new "2022-Q1" detail
#    sheet (same shape as the other quarterly sheets) and add a brand new
#    "总计" sheet after it with the 2022-Q1 row prepended to the history.

$wb = $excel.ActiveWorkbook

# A "clean" cell (default style, never touched) - pasting its format onto a
# cell strips any NumberFormat/quote-prefix residue while leaving the cell's
# value/type alone, so text-like numbers (e.g. "513690") keep their string
# type without round-tripping through a funny custom style.
$template = $wb.Worksheets.Item("2021-Q3")
$blank = $template.Range("Z100")

function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $blank.Copy()
    $range.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# 1. Repurpose the old "总计" sheet as the new "2022-Q1" detail sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Clear the old summary-table contents (B1:D4 header + 3 data rows) before
# laying out the quarterly fund-holding detail table.
$q1.Range("A1:D4").Clear()

# Reuse the header/index formatting (bold, bordered, centered) from an
# existing quarterly sheet.
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$q1.Range("A2").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0
Set-TextValue $q1.Range("B2") "513690"
Set-TextValue $q1.Range("C2") "博时恒生港股通高股息率ETF"
Set-TextValue $q1.Range("D2") "4.60"
Set-TextValue $q1.Range("E2") "99.64"
Set-TextValue $q1.Range("F2") "9.17"
Set-TextValue $q1.Range("G2") "0.4218"
$q1.Range("H2").Value = 1

# ---------------------------------------------------------------------
# 2. Brand new "总计" sheet (after "2022-Q1"), with the history table
#    incl. the new 2022-Q1 row on top.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$template.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$total.Range("A2:A5").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$rows = @(
    @(0, "2022-Q1", 1, 0.42),
    @(1, "2021-Q4", 3, 0.23),
    @(2, "2021-Q3", 2, 0.05),
    @(3, "2020-Q4", 6, 0.18)
)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 2 + $i
    $data = $rows[$i]
    $total.Range("A$r").Value = $data[0]
    Set-TextValue $total.Range("B$r") $data[1]
    $total.Range("C$r").Value = $data[2]
    $total.Range("D$r").Value = $data[3]
}

$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36
$total.Outline.SummaryRow = 1
$total.Outline.SummaryColumn = 1

# Keep the originally-active sheet selected (creating/renaming sheets above
# shifts the active tab onto whichever one we touched last).
$wb.Worksheets.Item("2020-Q4").Activate()

Write-Output "2022-Q1 data added"
